$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "isTargetFoe" (bool) column as column H
$ws.Range("H2").Value = "isTargetFoe"
$ws.Range("H3").Value = "bool"
$ws.Range("H4").Value = $true
$ws.Range("H5").Value = $true
$ws.Range("H6").Value = $true
$ws.Range("H7").Value = $true

# Match column width / outline metadata as closely as the host allows
$ws.Columns.Item(8).ColumnWidth = 11.71
$ws.Rows.Item(7).OutlineLevel = 6
$ws.Columns.Item(8).OutlineLevel = 7

# Match the recorded selection in the saved workbook
$ws.Range("H5").Select() | Out-Null
